$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) are stored as text in the source
# workbook (values like "26.216.22" or "  -1.46%  " are not valid Excel
# numbers/are intentionally kept as literal strings). Temporarily mark the
# range as Text so the COM layer does not auto-coerce numeric-looking
# values (e.g. "218.71") into real numbers, then restore the default
# "Normal" style afterwards so no stray per-cell formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '26.216.22'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '1.661.66'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '218.71'
$ws.Range("E5").Value = '  +0.68%  '
$ws.Range("D6").Value = '0.5230'
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("D7").Value = '1.006'
$ws.Range("E7").Value = '  +0.34%  '
$ws.Range("D8").Value = '0.2671'
$ws.Range("E8").Value = '  -0.35%  '
$ws.Range("D9").Value = '0.06321'
$ws.Range("E9").Value = '  -1.52%  '
$ws.Range("D10").Value = '21.08'
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("D11").Value = '0.07721'
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = '1.670.58'
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("D13").Value = '4.429'
$ws.Range("E13").Value = '  -1.72%  '
$ws.Range("D14").Value = '1.887.64'
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").Value = '0.5473'
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '0.0₅8216'
$ws.Range("E16").Value = '  -2.34%  '
$ws.Range("D17").Value = '64.89'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").Value = '26.246.57'
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("D20").Value = '4.666'
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("D21").Value = '193.13'
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("D22").Value = '10.16'
$ws.Range("E22").Value = '  -2.49%  '
$ws.Range("D23").Value = '6.089'
$ws.Range("E23").Value = '  -4.26%  '
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").Value = '138.76'
$ws.Range("E25").Value = '  -3.89%  '
$ws.Range("D26").Value = '0.1239'
$ws.Range("E26").Value = '  -3.03%  '
$ws.Range("D27").Value = '7.232'
$ws.Range("E27").Value = '  -3.21%  '
$ws.Range("D28").Value = '16.15'
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("D29").Value = '1.416'
$ws.Range("E29").Value = '  -1.42%  '
$ws.Range("D30").Value = '0.06006'
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("D31").Value = '1.282'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '3.652'
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("D33").Value = '3.309'
$ws.Range("E33").Value = '  -4.52%  '
$ws.Range("D34").Value = '1.634'
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("D35").Value = '0.9800'
$ws.Range("E35").Value = '  -3.13%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").Value = '2.786'
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.412'
$ws.Range("E37").Value = '  -0.48%  '
$ws.Range("D38").Value = '0.5872'
$ws.Range("E38").Value = '  +2.19%  '
$ws.Range("D39").Value = '0.01588'
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").Value = '5.949'
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("D41").Value = '0.8624'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D43").Value = '1.033.13'
$ws.Range("E43").Value = '  -3.88%  '
$ws.Range("D44").Value = '99.62'
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").Value = '1.801.88'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("E46").Value = '  +1.84%  '
$ws.Range("D47").Value = '57.13'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '8.102'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("D50").Value = '0.05183'
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '1.463'
$ws.Range("E51").Value = '  -0.42%  '

# Restore the original (unformatted) style now that the text values are in
# place, so the cells end up with no explicit style override - matching
# the source workbook.
$dataRange.Style = "Normal"
